# Update column F (dSF) values per the data repull / push-all-data / mean-calculation commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    4  = -2
    5  = -1
    6  = 1
    7  = -2
    8  = 1
    10 = 2
    11 = 4
    12 = 1
    13 = 1
    14 = 1
    15 = -3
    16 = -1
    17 = -1
    19 = 1
    20 = -4
    21 = -1
    22 = 5
    23 = 0
    24 = 3
    26 = -1
    27 = 2
    28 = 3
    29 = 8
    30 = 7
    31 = 8
    32 = 1
    34 = -4
    35 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
